# Leidingplanning 2017 - shift Saturday schedule header dates by 5 days
# (commit: "Better exception handling for missing dates and fix input schedule for Leiding")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
    $ws = $wb.Worksheets.Item("Leiding Zaterdag")
} catch {}

$ws.Range("C1").Value = 43029.3958333333
$ws.Range("D1").Value = 43029.4027777778
$ws.Range("E1").Value = 43029.4097222222
$ws.Range("F1").Value = 43029.4166666667
$ws.Range("G1").Value = 43029.4236111111
$ws.Range("H1").Value = 43029.4305555556
$ws.Range("I1").Value = 43029.4375
$ws.Range("J1").Value = 43029.4444444445
$ws.Range("K1").Value = 43029.4513888889
$ws.Range("L1").Value = 43029.4583333333
$ws.Range("M1").Value = 43029.4652777778
$ws.Range("N1").Value = 43029.4722222222
$ws.Range("O1").Value = 43029.4791666667
$ws.Range("P1").Value = 43029.4861111111
$ws.Range("Q1").Value = 43029.4930555556
$ws.Range("R1").Value = 43029.5
$ws.Range("S1").Value = 43029.5069444445
$ws.Range("T1").Value = 43029.5138888889
$ws.Range("U1").Value = 43029.5208333333
$ws.Range("V1").Value = 43029.5277777778
$ws.Range("W1").Value = 43029.5347222222
$ws.Range("X1").Value = 43029.5416666667
$ws.Range("Y1").Value = 43029.5486111111
$ws.Range("Z1").Value = 43029.5555555556
$ws.Range("AA1").Value = 43029.5625
$ws.Range("AB1").Value = 43029.5694444445
$ws.Range("AC1").Value = 43029.5763888889
$ws.Range("AD1").Value = 43029.5833333334
$ws.Range("AE1").Value = 43029.5902777778
$ws.Range("AF1").Value = 43029.5972222222
$ws.Range("AG1").Value = 43029.6041666667
$ws.Range("AH1").Value = 43029.6111111111
$ws.Range("AI1").Value = 43029.6180555556
$ws.Range("AJ1").Value = 43029.625
$ws.Range("AK1").Value = 43029.6319444445
$ws.Range("AL1").Value = 43029.6388888889
$ws.Range("AM1").Value = 43029.6458333334
$ws.Range("AN1").Value = 43029.6527777778
$ws.Range("AO1").Value = 43029.6597222223
$ws.Range("AP1").Value = 43029.6666666667
$ws.Range("AQ1").Value = 43029.6736111111
$ws.Range("AR1").Value = 43029.6805555556
$ws.Range("AS1").Value = 43029.6875
$ws.Range("AT1").Value = 43029.6944444445
$ws.Range("AU1").Value = 43029.7013888889
$ws.Range("AV1").Value = 43029.7083333334
$ws.Range("AW1").Value = 43029.7152777778
$ws.Range("AX1").Value = 43029.7222222223
$ws.Range("AY1").Value = 43029.7291666667
$ws.Range("AZ1").Value = 43029.7361111111
$ws.Range("BA1").Value = 43029.7430555556
$ws.Range("BB1").Value = 43029.75
$ws.Range("BC1").Value = 43029.7569444445
$ws.Range("BD1").Value = 43029.7638888889
$ws.Range("BE1").Value = 43029.7708333334
$ws.Range("BF1").Value = 43029.7777777778
$ws.Range("BG1").Value = 43029.7847222223
$ws.Range("BH1").Value = 43029.7916666667
$ws.Range("BI1").Value = 43029.7986111111
$ws.Range("BJ1").Value = 43029.8055555556
$ws.Range("BK1").Value = 43029.8125
$ws.Range("BL1").Value = 43029.8194444445
$ws.Range("BM1").Value = 43029.8263888889
$ws.Range("BN1").Value = 43029.8333333334
$ws.Range("BO1").Value = 43029.8402777778
$ws.Range("BP1").Value = 43029.8472222223
$ws.Range("BQ1").Value = 43029.8541666667
$ws.Range("BR1").Value = 43029.8611111112
$ws.Range("BS1").Value = 43029.8680555556
$ws.Range("BT1").Value = 43029.875
$ws.Range("BU1").Value = 43029.8819444445
$ws.Range("BV1").Value = 43029.8888888889
$ws.Range("BW1").Value = 43029.8958333334
$ws.Range("BX1").Value = 43029.9027777778
$ws.Range("BY1").Value = 43029.9097222223
$ws.Range("BZ1").Value = 43029.9166666667
$ws.Range("CA1").Value = 43029.9236111112
$ws.Range("CB1").Value = 43029.9305555556
$ws.Range("CC1").Value = 43029.9375
$ws.Range("CD1").Value = 43029.9444444445
$ws.Range("CE1").Value = 43029.9513888889
$ws.Range("CF1").Value = 43029.9583333334
$ws.Range("CG1").Value = 43029.9652777778
$ws.Range("CH1").Value = 43029.9722222223
$ws.Range("CI1").Value = 43029.9791666667
$ws.Range("CJ1").Value = 43029.9861111112
$ws.Range("CK1").Value = 43029.9930555556
$ws.Range("CL1").Value = 43030.0000000001

# Update the view to match where the edit left the selection/scroll position
try {
    $ws.Range("CK1").Select() | Out-Null
} catch {}
try {
    $excel.ActiveWindow.ScrollColumn = 82
} catch {}
